$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H70").Value = 4204.15
$ws.Range("I70").Value = 4148.857
$ws.Range("J70").Value = 4333.1665
$ws.Range("K70").Value = 12446.571
$ws.Range("L70").Value = 12999.4995
$ws.Range("M70").Value = -12176.571
$ws.Range("N70").Value = -13539.4995
$ws.Range("H73").Value = 4204.15
$ws.Range("I73").Value = 4148.857
$ws.Range("J73").Value = 4333.1665
$ws.Range("K73").Value = 12446.571
$ws.Range("L73").Value = 12999.4995
$ws.Range("M73").Value = -11510.571
$ws.Range("N73").Value = -14871.4995
$ws.Range("H132").Value = 4380.778
$ws.Range("I132").Value = 4547.5884
$ws.Range("K132").Value = 13642.7652
$ws.Range("M132").Value = -11112.7652
$ws.Range("H135").Value = 1523.7646
$ws.Range("I135").Value = 1594
$ws.Range("K135").Value = 14346
$ws.Range("M135").Value = -11811

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4688.8
$ws.Range("I61").Value = 0
$ws.Range("J61").Value = 4688.8
$ws.Range("K61").Value = 0
$ws.Range("L61").Value = 4688.8
$ws.Range("M61").ClearContents()
$ws.Range("N61").Value = -5112.8
$ws.Range("H74").Value = 181131.28
$ws.Range("I74").Value = 209819.83
$ws.Range("J74").Value = 9000
$ws.Range("K74").Value = 209819.83
$ws.Range("L74").Value = 9000
$ws.Range("M74").Value = -208945.83
$ws.Range("N74").Value = -10748
$ws.Range("H77").Value = 181131.28
$ws.Range("I77").Value = 209819.83
$ws.Range("J77").Value = 9000
$ws.Range("K77").Value = 1049099.15
$ws.Range("L77").Value = 45000
$ws.Range("M77").Value = -1044731.15
$ws.Range("N77").Value = -53736
$ws.Range("H114").Value = 0
$ws.Range("J114").Value = 0
$ws.Range("L114").Value = 0
$ws.Range("N114").ClearContents()
$ws.Range("H132").Value = 27924.656
$ws.Range("I132").Value = 6214.161
$ws.Range("J132").Value = 700950
$ws.Range("K132").Value = 18642.483
$ws.Range("L132").Value = 2102850
$ws.Range("M132").Value = -16112.483
$ws.Range("N132").Value = -2107910
$ws.Range("H133").Value = 107000
$ws.Range("J133").Value = 107000
$ws.Range("L133").Value = 107000
$ws.Range("N133").Value = -112060
$ws.Range("H136").Value = 4688.8
$ws.Range("I136").Value = 0
$ws.Range("J136").Value = 4688.8
$ws.Range("K136").Value = 0
$ws.Range("L136").Value = 14066.4
$ws.Range("M136").ClearContents()
$ws.Range("N136").Value = -19166.4

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 1687.1111
$ws.Range("I134").Value = 1528.5
$ws.Range("K134").Value = 4585.5
$ws.Range("M134").Value = -2050.5
$ws.Range("H141").Value = 18354.5
$ws.Range("I141").Value = 18354.5
$ws.Range("K141").Value = 18354.5
$ws.Range("M141").Value = -13174.5

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 131273000
$ws.Range("I4").Value = 7169143
$ws.Range("K4").Value = 7169143
$ws.Range("M4").Value = -7169031
$ws.Range("H41").Value = 28881.125
$ws.Range("I41").Value = 4529.5
$ws.Range("J41").Value = 36998.332
$ws.Range("K41").Value = 4529.5
$ws.Range("L41").Value = 36998.332
$ws.Range("M41").Value = -4101.5
$ws.Range("N41").Value = -37854.332
$ws.Range("H106").Value = 29360.334
$ws.Range("I106").Value = 34610
$ws.Range("J106").Value = 26735.5
$ws.Range("K106").Value = 34610
$ws.Range("L106").Value = 26735.5
$ws.Range("M106").Value = -33348
$ws.Range("N106").Value = -29259.5
$ws.Range("H131").Value = 29449.75
$ws.Range("J131").Value = 31942.572
$ws.Range("L131").Value = 31942.572
$ws.Range("N131").Value = -42022.572
$ws.Range("H134").Value = 2532.175
$ws.Range("I134").Value = 2220.2856
$ws.Range("K134").Value = 6660.8568
$ws.Range("M134").Value = -4125.8568
$ws.Range("H141").Value = 104392.09
$ws.Range("I141").Value = 35000
$ws.Range("J141").Value = 130414.125
$ws.Range("K141").Value = 35000
$ws.Range("L141").Value = 130414.125
$ws.Range("M141").Value = -29820
$ws.Range("N141").Value = -140774.125

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 41609936
$ws.Range("I4").Value = 75367700
$ws.Range("K4").Value = 226103100
$ws.Range("M4").Value = -226102988
$ws.Range("H122").Value = 1046.5358
$ws.Range("J122").Value = 1157
$ws.Range("L122").Value = 10413
$ws.Range("N122").Value = -15313

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 641203.8
$ws.Range("J33").Value = 1038673
$ws.Range("L33").Value = 1038673
$ws.Range("N33").Value = -1039177
$ws.Range("H70").Value = 85489.56
$ws.Range("I70").Value = 115646.39
$ws.Range("K70").Value = 115646.39
$ws.Range("M70").Value = -115376.39
$ws.Range("H73").Value = 85489.56
$ws.Range("I73").Value = 115646.39
$ws.Range("K73").Value = 115646.39
$ws.Range("M73").Value = -114710.39
$ws.Range("H80").Value = 2259.875
$ws.Range("I80").Value = 1775.8
$ws.Range("J80").Value = 3066.6667
$ws.Range("K80").Value = 1775.8
$ws.Range("L80").Value = 3066.6667
$ws.Range("M80").Value = -777.8
$ws.Range("N80").Value = -5062.6667
$ws.Range("H83").Value = 2259.875
$ws.Range("I83").Value = 1775.8
$ws.Range("J83").Value = 3066.6667
$ws.Range("K83").Value = 8879
$ws.Range("L83").Value = 15333.3335
$ws.Range("M83").Value = -3887
$ws.Range("N83").Value = -25317.3335
$ws.Range("H102").Value = 6170.727
$ws.Range("I102").Value = 2633.5715
$ws.Range("K102").Value = 2633.5715
$ws.Range("M102").Value = -1011.5715
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("L132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("N132").ClearContents()

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6776.1665
$ws.Range("I7").Value = 6254.3335
$ws.Range("K7").Value = 6254.3335
$ws.Range("M7").Value = -6142.3335
$ws.Range("H61").Value = 26318490
$ws.Range("I61").Value = 29414472
$ws.Range("J61").Value = 2650
$ws.Range("K61").Value = 29414472
$ws.Range("L61").Value = 2650
$ws.Range("M61").Value = -29414270
$ws.Range("N61").Value = -3054
$ws.Range("H64").Value = 50562.5
$ws.Range("J64").Value = 50562.5
$ws.Range("L64").Value = 50562.5
$ws.Range("N64").Value = -51012.5
$ws.Range("H67").Value = 50562.5
$ws.Range("J67").Value = 50562.5
$ws.Range("L67").Value = 50562.5
$ws.Range("N67").Value = -52122.5
$ws.Range("H68").Value = 2450.1667
$ws.Range("I68").Value = 2500.3333
$ws.Range("J68").Value = 2400
$ws.Range("K68").Value = 2500.3333
$ws.Range("L68").Value = 2400
$ws.Range("M68").Value = -1751.3333
$ws.Range("N68").Value = -3898
$ws.Range("H71").Value = 2450.1667
$ws.Range("I71").Value = 2500.3333
$ws.Range("J71").Value = 2400
$ws.Range("K71").Value = 12501.6665
$ws.Range("L71").Value = 12000
$ws.Range("M71").Value = -8757.666499999999
$ws.Range("N71").Value = -19488
$ws.Range("H82").Value = 4122
$ws.Range("I82").Value = 2369
$ws.Range("J82").Value = 5875
$ws.Range("K82").Value = 2369
$ws.Range("L82").Value = 5875
$ws.Range("M82").Value = -2008
$ws.Range("N82").Value = -6597
$ws.Range("H85").Value = 4122
$ws.Range("I85").Value = 2369
$ws.Range("J85").Value = 5875
$ws.Range("K85").Value = 2369
$ws.Range("L85").Value = 5875
$ws.Range("M85").Value = -1121
$ws.Range("N85").Value = -8371
$ws.Range("H100").Value = 3436.6191
$ws.Range("J100").Value = 3899.8572
$ws.Range("L100").Value = 3899.8572
$ws.Range("N100").Value = -4981.8572
$ws.Range("H113").Value = 26318490
$ws.Range("I113").Value = 29414472
$ws.Range("J113").Value = 2650
$ws.Range("K113").Value = 29414472
$ws.Range("L113").Value = 2650
$ws.Range("M113").Value = -29412302
$ws.Range("N113").Value = -6990
$ws.Range("H122").Value = 2810.8484
$ws.Range("I122").Value = 2602.48
$ws.Range("K122").Value = 7807.440000000001
$ws.Range("M122").Value = -5357.440000000001
$ws.Range("H126").Value = 6776.1665
$ws.Range("I126").Value = 6254.3335
$ws.Range("K126").Value = 18763.0005
$ws.Range("M126").Value = -16293.0005
$ws.Range("H136").Value = 1952.0286
$ws.Range("I136").Value = 1720.0646
$ws.Range("J136").Value = 3749.75
$ws.Range("K136").Value = 5160.1938
$ws.Range("L136").Value = 11249.25
$ws.Range("M136").Value = -2610.1938
$ws.Range("N136").Value = -16349.25

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 6202.385
$ws.Range("I62").Value = 5126.5
$ws.Range("J62").Value = 7124.5713
$ws.Range("K62").Value = 5126.5
$ws.Range("L62").Value = 7124.5713
$ws.Range("M62").Value = -4502.5
$ws.Range("N62").Value = -8372.5713
$ws.Range("H65").Value = 6202.385
$ws.Range("I65").Value = 5126.5
$ws.Range("J65").Value = 7124.5713
$ws.Range("K65").Value = 25632.5
$ws.Range("L65").Value = 35622.85649999999
$ws.Range("M65").Value = -22512.5
$ws.Range("N65").Value = -41862.85649999999
$ws.Range("H81").Value = 1388.05
$ws.Range("I81").Value = 1391.6
$ws.Range("J81").Value = 1377.4
$ws.Range("K81").Value = 2783.2
$ws.Range("L81").Value = 2754.8
$ws.Range("M81").Value = -1722.2
$ws.Range("N81").Value = -4876.8
$ws.Range("H84").Value = 1388.05
$ws.Range("I84").Value = 1391.6
$ws.Range("J84").Value = 1377.4
$ws.Range("K84").Value = 13916
$ws.Range("L84").Value = 13774
$ws.Range("M84").Value = -8612
$ws.Range("N84").Value = -24382
$ws.Range("H122").Value = 5636.826
$ws.Range("I122").Value = 5897.6
$ws.Range("J122").Value = 3898.3333
$ws.Range("K122").Value = 17692.8
$ws.Range("L122").Value = 11694.9999
$ws.Range("M122").Value = -15242.8
$ws.Range("N122").Value = -16594.9999
$ws.Range("H136").Value = 4611.1763
$ws.Range("I136").Value = 3125.88
$ws.Range("K136").Value = 9377.639999999999
$ws.Range("M136").Value = -6827.639999999999
